$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin -> Bitcoin
$ws.Range("D2").Value = "28.200.21"
$ws.Range("E2").Value = "  -0.99%  "

# Row 3: Ethereum -> Ethereum
$ws.Range("D3").Value = "1.796.73"
$ws.Range("E3").Value = "  -1.33%  "

# Row 4: TetherUSD -> TetherUSD
$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").Style = "Normal"

# Row 5: BNB -> BNB
$ws.Range("D5").Value = "'314.64"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.33%  "

# Row 6: USDC -> USDC
$ws.Range("D6").Value = "'1.003"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.06%  "

# Row 7: XRP -> XRP
$ws.Range("D7").Value = "'0.5205"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.16%  "

# Row 8: Cardano -> Cardano
$ws.Range("D8").Value = "'0.3814"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.07%  "

# Row 9: Dogecoin -> Dogecoin
$ws.Range("D9").Value = "'0.07902"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.82%  "

# Row 10: OKB -> OKB
$ws.Range("D10").Value = "'41.56"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.01%  "

# Row 11: Polygon -> Polygon
$ws.Range("D11").Value = "'1.096"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.01%  "

# Row 12: Polkadot -> Polkadot
$ws.Range("D12").Value = "'6.267"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.17%  "

# Row 13: BinanceUSD -> BinanceUSD
$ws.Range("D13").Value = "'1.003"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.12%  "

# Row 14: Solana -> Solana
$ws.Range("D14").Value = "'20.57"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.05%  "

# Row 15: Chainlink -> WrappedEther
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.792.50"
$ws.Range("E15").Value = "  -2.04%  "

# Row 16: WrappedEther -> Chainlink
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'7.273"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.09%  "

# Row 17: Litecoin -> Litecoin
$ws.Range("D17").Value = "'93.42"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.51%  "

# Row 18: ShibaInu -> ShibaInu
$ws.Range("D18").Value = "'0.00001081"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.99%  "

# Row 19: TRON -> TRON
$ws.Range("D19").Value = "'0.06552"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.19%  "

# Row 20: Dai -> Dai
$ws.Range("D20").Value = "'1.003"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.05%  "

# Row 21: Avalanche -> Avalanche
$ws.Range("D21").Value = "'17.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.39%  "

# Row 22: Uniswap -> Uniswap
$ws.Range("D22").Value = "'5.961"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.07%  "

# Row 23: WrappedBTC -> WrappedBTC
$ws.Range("D23").Value = "28.262.59"
$ws.Range("E23").Value = "  -0.91%  "

# Row 24: Cosmos -> Cosmos
$ws.Range("D24").Value = "'11.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.95%  "

# Row 25: Toncoin -> Toncoin
$ws.Range("D25").Value = "'2.264"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.42%  "

# Row 26: Monero -> Monero
$ws.Range("D26").Value = "'159.84"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.51%  "

# Row 27: EthereumClassic -> EthereumClassic
$ws.Range("D27").Value = "'20.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.49%  "

# Row 28: WrappedliquidstakedEther2.0 -> WrappedliquidstakedEther2.0
$ws.Range("D28").Value = "1.998.62"
$ws.Range("E28").Value = "  -1.91%  "

# Row 29: LidoDAOToken -> LidoDAOToken
$ws.Range("D29").Value = "'2.340"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.65%  "

# Row 30: BitcoinCash -> BitcoinCash
$ws.Range("D30").Value = "'123.44"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.15%  "

# Row 31: Stellar -> Stellar
$ws.Range("D31").Value = "'0.1074"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.42%  "

# Row 32: ImmutableX -> ImmutableX
$ws.Range("D32").Value = "'1.052"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.87%  "

# Row 33: HuobiToken -> HuobiToken
$ws.Range("D33").Value = "'3.678"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.33%  "

# Row 34: Filecoin -> Filecoin
$ws.Range("D34").Value = "'5.544"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.33%  "

# Row 35: Hedera -> Hedera
$ws.Range("D35").Value = "'0.07138"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.35%  "

# Row 36: Aptos -> VeChain
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").Value = "'0.02328"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.81%  "

# Row 37: VeChain -> Aptos
$ws.Range("B37").Value = "Aptos"
$ws.Range("C37").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D37").Value = "'11.85"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.54%  "

# Row 38: Algorand -> Algorand
$ws.Range("D38").Value = "'0.2140"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.49%  "

# Row 39: InternetComputer(DFINITY) -> FraxShare
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "'8.730"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.65%  "

# Row 40: FraxShare -> InternetComputer(DFINITY)
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").Value = "'5.073"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.71%  "

# Row 41: TheSandbox -> TheSandbox
$ws.Range("D41").Value = "'0.6209"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.93%  "

# Row 42: TrustWalletToken -> TrustWalletToken
$ws.Range("D42").Value = "'1.158"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.33%  "

# Row 43: WEMIXTOKEN -> WEMIXTOKEN
$ws.Range("D43").Value = "'1.369"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.49%  "

# Row 44: EnergySwap -> Decentraland
$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").Value = "'0.6127"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.45%  "

# Row 45: Decentraland -> EnergySwap
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'13.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.32%  "

# Row 46: PancakeSwap -> PancakeSwap
$ws.Range("D46").Value = "'3.774"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.10%  "

# Row 47: Quant -> Quant
$ws.Range("D47").Value = "'127.10"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.35%  "

# Row 48: EOS -> EOS
$ws.Range("D48").Value = "'1.228"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.75%  "

# Row 49: NEARProtocol -> NEARProtocol
$ws.Range("D49").Value = "'1.919"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.47%  "

# Row 50: Cronos -> Cronos
$ws.Range("D50").Value = "'0.06763"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.52%  "

# Row 51: Aave -> ThetaToken
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").Value = "'1.061"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.41%  "
